$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append new reference-table rows for the new "incident_asset_type" and
# "received_via" lookup categories (rows 51-56), following the existing
# pattern: column A = category key, column B = sequence number, column C = value.

$data = @(
    @("incident_asset_type", 1, "Hilang"),
    @("incident_asset_type", 2, "Rosak"),
    @("received_via", 1, "Telefon"),
    @("received_via", 2, "Emel"),
    @("received_via", 3, "Chatbot"),
    @("received_via", 4, "Live Chat")
)

$row = 51
foreach ($item in $data) {
    $ws.Cells.Item($row, 1).Value = $item[0]
    $ws.Cells.Item($row, 2).Value = $item[1]
    $ws.Cells.Item($row, 3).Value = $item[2]
    $row++
}

# Update the view to match: scrolled down with the last new cell selected.
$ws.Range("A52").Select()
